$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Joins($text) {
    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    return $text
}

$cellsToFix = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cellsToFix) {
    $rng = $ws.Range($addr)
    $rng.Value = Fix-Joins $rng.Value2
}

# Update sheet view: scroll so A6 is the top-left visible cell, and select C7
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 6

# Update column C width (no longer "best fit"; now a fixed custom width)
$ws.Columns.Item(3).ColumnWidth = 67.33
